# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-RowValues($row, $name, $values) {
    if ($name -ne $null) {
        $ws.Cells.Item($row, 1).Value = $name
    }
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 13:05"

# --- Catar (row 23): refreshed case counts ---
Set-RowValues 23 $null @(48947, 1740, 13283, 35634, 0, 2, 30)

# --- Refreshed case counts for a few other countries (rows 54, 64, 76) ---
$ws.Cells.Item(54, 4).Value = 4669
$ws.Cells.Item(54, 5).Value = 4598

$ws.Cells.Item(64, 4).Value = 6566
$ws.Cells.Item(64, 5).Value = 470

$ws.Cells.Item(76, 4).Value = 2659
$ws.Cells.Item(76, 5).Value = 660

# --- Nepal inserted ahead of Paraguay; Paraguay, Burkina Faso, Sudan del
#     Sur and Uruguay each shift down one row (rows 118-122) ---
Set-RowValues 118 "Nepal"         @(886, 114, 183, 699, 0, 0, 4)
Set-RowValues 119 "Paraguay"      @(877, 0,   382, 484, 0, 0, 11)
Set-RowValues 120 "Burkina Faso"  @(845, 13,  672, 120, 0, 1, 53)
Set-RowValues 121 "Sudan del Sur" @(806, 0,   6,   792, 0, 0, 8)
Set-RowValues 122 "Uruguay"       @(789, 0,   638, 129, 0, 0, 22)

# --- Madagascar / Malta swap order (rows 133-134), refreshed counts ---
Set-RowValues 133 "Madagascar" @(612, 26, 151, 459, 0, 0, 2)
Set-RowValues 134 "Malta"      @(612, 1,  491, 114, 0, 1, 7)
